# WEEK 05-18 last edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week")

# --- Actual hours (column I) updated for rows 7-11 ---
$ws.Range("I7").Value = 43
$ws.Range("I8").Value = 83
$ws.Range("I9").Value = 43
$ws.Range("I10").Value = 83
$ws.Range("I11").Value = 83

# --- Nudge the week-selector spin button back onto the pixel grid ---
# (from col D/row3 offset 12.0pt, to col E/row5 offset 6.0pt x 1.5pt)
$shp = $ws.Shapes.Item("weekSelector")
$shp.Top = 40.55
$shp.Left = 244.21875
$shp.Width = 14.34375
$shp.Height = 18.05

# --- View state: scroll position + active selection ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("K13").Select()
